# Updated cryptos list - apply Price (D) and Volume(1h) (E) changes per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "63.934.55"
$ws.Cells.Item(2, 5).Value = "  +0.44%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.637.34"
$ws.Cells.Item(3, 5).Value = "  +1.00%  "

$ws.Cells.Item(4, 5).Value = "  +0.12%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "580.00"
$ws.Cells.Item(5, 5).Value = "  +1.08%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "157.02"
$ws.Cells.Item(6, 5).Value = "  +1.64%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.632"
$ws.Cells.Item(7, 5).Value = "  +0.46%  "

$ws.Cells.Item(8, 5).Value = "  +0.09%  "

$ws.Cells.Item(9, 5).Value = "  -0.96%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.386"
$ws.Cells.Item(11, 5).Value = "  +0.47%  "

$ws.Cells.Item(12, 5).Value = "  +1.03%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "28.84"
$ws.Cells.Item(13, 5).Value = "  +2.56%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "3.113.07"
$ws.Cells.Item(14, 5).Value = "  +1.04%  "

$ws.Cells.Item(15, 5).Value = "  +1.98%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "63.854.33"
$ws.Cells.Item(16, 5).Value = "  +0.50%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.643.25"
$ws.Cells.Item(17, 5).Value = "  +0.87%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "12.21"
$ws.Cells.Item(18, 5).Value = "  +1.45%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.77"
$ws.Cells.Item(19, 5).Value = "  +4.13%  "

$ws.Cells.Item(20, 5).Value = "  -1.69%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "344.62"
$ws.Cells.Item(21, 5).Value = "  +0.52%  "

$ws.Cells.Item(22, 5).Value = "  +0.12%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "68.43"
$ws.Cells.Item(23, 5).Value = "  +2.27%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.90"
$ws.Cells.Item(24, 5).Value = "  +9.31%  "

$ws.Cells.Item(25, 5).Value = "  +5.64%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.65"
$ws.Cells.Item(26, 5).Value = "  +4.64%  "

$ws.Cells.Item(27, 5).Value = "  +0.58%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "584.74"
$ws.Cells.Item(28, 5).Value = "  +1.43%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "8.35"
$ws.Cells.Item(29, 5).Value = "  +5.87%  "

$ws.Cells.Item(30, 5).Value = "  +0.65%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.977"
$ws.Cells.Item(31, 5).Value = "  -2.38%  "

$ws.Cells.Item(32, 5).Value = "  +0.22%  "

$ws.Cells.Item(33, 5).Value = "  +3.32%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "6.69"
$ws.Cells.Item(34, 5).Value = "  +3.57%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.50"
$ws.Cells.Item(35, 5).Value = "  +3.73%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.405"
$ws.Cells.Item(36, 5).Value = "  -0.79%  "

$ws.Cells.Item(37, 5).Value = "  -0.06%  "

$ws.Cells.Item(38, 5).Value = "  +0.11%  "

$ws.Cells.Item(39, 5).Value = "  +3.69%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "154.46"
$ws.Cells.Item(40, 5).Value = "  +0.53%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.59"
$ws.Cells.Item(41, 5).Value = "  +9.68%  "

$ws.Cells.Item(42, 5).Value = "  +0.03%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "162.90"
$ws.Cells.Item(43, 5).Value = "  +4.38%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "24.29"
$ws.Cells.Item(44, 5).Value = "  +6.30%  "

$ws.Cells.Item(45, 5).Value = "  -0.26%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0590"
$ws.Cells.Item(46, 5).Value = "  -0.44%  "

$ws.Cells.Item(47, 5).Value = "  +1.11%  "

$ws.Cells.Item(48, 5).Value = "  -0.48%  "

$ws.Cells.Item(49, 5).Value = "  +0.04%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0₆0239"
$ws.Cells.Item(50, 5).Value = "  +3.12%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.793"
$ws.Cells.Item(51, 5).Value = "  +3.25%  "
